$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the latest week (2022-01-25, serial 44586),
# pushing all existing data rows down by 2.
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# Row 15: Black Amber, Primera
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44586
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103002
$ws.Range("J15").Value = "Ciruela"
$ws.Range("K15").Value = "Black Amber"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 9000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 9500
$ws.Range("Q15").Value = "$/bandeja 18 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 528
$ws.Range("T15").Value = 18

# Row 16: Black Amber, Segunda
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44586
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100103
$ws.Range("H16").Value = "Frutos de hueso (carozo)"
$ws.Range("I16").Value = 100103002
$ws.Range("J16").Value = "Ciruela"
$ws.Range("K16").Value = "Black Amber"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("Q16").Value = "$/bandeja 18 kilos granel"
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 444
$ws.Range("T16").Value = 18
